$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: replace the 12-column header row with a single title cell ---
$ws.Range("A1:L1").ClearContents()
$ws.Range("A1").Value = "Herald College Kathmandu"

# --- Data rows 2-10: reorder / drop columns ---
# Old layout: A Day, B Time, C Code, D Title, E Hours, F Type, G Lecturer,
#             H Room, I Block, J Group, K Level, L Course
# New layout: A Day, B Time, C Hours, D Code, E Title, F Type, G Lecturer,
#             H Group, I Block, J Room
# (old K Level and L Course columns are dropped)
for ($r = 2; $r -le 10; $r++) {
    $code  = $ws.Cells.Item($r, 3).Value2   # C: Module Code
    $title = $ws.Cells.Item($r, 4).Value2   # D: Module Title
    $hours = $ws.Cells.Item($r, 5).Value2   # E: Hours
    $room  = $ws.Cells.Item($r, 8).Value2   # H: Room
    $group = $ws.Cells.Item($r, 10).Value2  # J: Group

    $ws.Cells.Item($r, 3).Value  = $hours   # C: Hours
    $ws.Cells.Item($r, 4).Value  = $code    # D: Module Code
    $ws.Cells.Item($r, 5).Value  = $title   # E: Module Title
    # F (Class Type) and G (Lecturer) stay in place
    $ws.Cells.Item($r, 8).Value  = $group   # H: Group
    # I (Block) stays in place
    $ws.Cells.Item($r, 10).Value = $room    # J: Room

    # Drop old K (Level) and L (Course) columns
    $ws.Cells.Item($r, 11).ClearContents()
    $ws.Cells.Item($r, 12).ClearContents()
}

